# Generate Report for Handoff
# Updates the localization-status workbook: marks b.md as "Ready for handoff"
# with a freshly generated handoff xliff + timestamp for both zh-cn and
# de-de locales, and records the "stale handback" error detail message.

$wb = $excel.ActiveWorkbook

$readyStatus = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/b5fd1b853ad62aaa8b165ad6d2871363ce0234a3/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/c167867fe3c33551f761c77c9ba6902ba69d8e6a/e2e/b.md."

# --- Sheet "Overview" (b.md row) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $readyStatus
$wsOverview.Range("F3").Value = $readyStatus
$wsOverview.Range("G3").Value = "2016-08-13 06:45:24"

# --- Sheet "zh-cn" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# a.md row stays handed off -> now "Ready for handoff"
$wsZhCn.Range("C2").Value = $readyStatus

# b.md row: new handoff file + timestamp, stale handback error detail
$wsZhCn.Range("C3").Value = $readyStatus
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-13 06:45:17"
$wsZhCn.Range("P3").Value = $errorDetail

# --- Sheet "de-de" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17

# b.md row: new handoff file + timestamp, stale handback error detail
$wsDeDe.Range("C3").Value = $readyStatus
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-13 06:45:24"
$wsDeDe.Range("P3").Value = $errorDetail
